$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Overfitting in Automated Program Repair..." entry) is being re-annotated:
# the previously recorded bias text is replaced with the correct one, and the
# "omitted bias?" column (E) is now filled in with "No".
$ws.Range("D16").Value = "It reported the ""Only-manual validation bias"" and ""Only-independent test validation bias"", the same biases reported by Le et al. \cite{le2019reliability}."
$ws.Range("E16").Value = "No"

# Reflect where the author was working when making the edit: scrolled down to
# row 7 and had D12 selected.
$ws.Range("D12").Select() | Out-Null
